# Update "Lương" sheet: remove the old column B values and insert new
# rows for "Ứng lương" per location plus grand-total rows at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Clear the whole sheet first, then rewrite column A from scratch so the
# row layout matches exactly (inserted rows shift everything below them).
$ws.Cells.Clear()

$labels = @(
    "Danh mục",
    "Ngày công",
    "Phụ cấp",
    "Lương cơ bản tại CẦN THƠ",
    "Chiết khấu sale chính tại CẦN THƠ",
    "Chiết khấu sale phụ tại CẦN THƠ",
    "Đơn 1 bác sĩ tại CẦN THƠ",
    "Đơn 2 bác sĩ tại CẦN THƠ",
    "Công phụ phẫu 1 tại CẦN THƠ",
    "Công phụ phẫu 2 tại CẦN THƠ",
    "Ứng lương tại CẦN THƠ",
    "Lương cơ bản tại LONG XUYÊN",
    "Chiết khấu sale chính tại LONG XUYÊN",
    "Chiết khấu sale phụ tại LONG XUYÊN",
    "Đơn 1 bác sĩ tại LONG XUYÊN",
    "Đơn 2 bác sĩ tại LONG XUYÊN",
    "Công phụ phẫu 1 tại LONG XUYÊN",
    "Công phụ phẫu 2 tại LONG XUYÊN",
    "Ứng lương tại LONG XUYÊN",
    "Lương cơ bản tại SÓC TRĂNG",
    "Chiết khấu sale chính tại SÓC TRĂNG",
    "Chiết khấu sale phụ tại SÓC TRĂNG",
    "Đơn 1 bác sĩ tại SÓC TRĂNG",
    "Đơn 2 bác sĩ tại SÓC TRĂNG",
    "Công phụ phẫu 1 tại SÓC TRĂNG",
    "Công phụ phẫu 2 tại SÓC TRĂNG",
    "Ứng lương tại SÓC TRĂNG",
    "Tổng lương tại CẦN THƠ",
    "Tổng lương tại LONG XUYÊN",
    "Tổng lương tại SÓC TRĂNG",
    "Tổng lương"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $labels[$i]
}

$ws.Range("A1").Select()
